$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# ---------------------------------------------------------------------------
# 1. New header cells in row 1 (columns V..AI, with blanks at U/AA/AG)
# ---------------------------------------------------------------------------
$ws.Range("V1").Value  = "ending HRU volume (m3)"
$ws.Range("W1").Value  = "ending reach volume (m3)"
$ws.Range("X1").Value  = "ending reservoir volume (m3)"
$ws.Range("Y1").Value  = "ending volume (m3)"
$ws.Range("Z1").Value  = "ending volume (mm)"
$ws.Range("AB1").Value = "starting HRU volume (m3)"
$ws.Range("AC1").Value = "starting reach volume (m3)"
$ws.Range("AD1").Value = "starting reservoir volume (m3)"
$ws.Range("AE1").Value = "starting volume (m3)"
$ws.Range("AF1").Value = "starting volume (mm)"
$ws.Range("AH1").Value = "change in volume (m3)"
$ws.Range("AI1").Value = "change in volume (mm)"

# ---------------------------------------------------------------------------
# 2. New column widths
# ---------------------------------------------------------------------------
$ws.Columns("V").ColumnWidth = 11
$ws.Columns("X").ColumnWidth = 11.88671875
$ws.Columns("Y").ColumnWidth = 11

# ---------------------------------------------------------------------------
# 3. Insert three blank rows above the old row 125 (pushes old 125->128,
#    126->129, 127->130, 128->131) and restore the blank rows' formatting
#    to match the plain "spacer row" style used elsewhere in the sheet.
# ---------------------------------------------------------------------------
$ws.Rows("125:127").Insert()

foreach ($r in 125..127) {
    foreach ($col in @("D","E","F","G","H","I","J","K","L","M","N","Q")) {
        $ws.Range($col + $r).NumberFormat = "0.00"
    }
    $ws.Range("O" + $r).NumberFormat = "0"
    $ws.Range("P" + $r).NumberFormat = "0"
    $ws.Range("R" + $r).NumberFormat = "0.000000"
    $ws.Range("B" + $r).WrapText = $false
}

# ---------------------------------------------------------------------------
# 4. Populate the new row 125 with the "CW3M C794" model run data
# ---------------------------------------------------------------------------
$ws.Range("A125").Value = "???"
$ws.Range("B125").Value = "Demo_Baseline 2010-18"
$ws.Range("B125").WrapText = $true
$ws.Range("C125").Value = "2010-18"

$ws.Range("D125").Value = 946.24761266666656
$ws.Range("E125").Value = 1890.2624918888889
$ws.Range("F125").Value = 0.56412044444444454
$ws.Range("G125").Value = 270.41205844444437
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0.20218555555555556
$ws.Range("J125").Value = 8.1970344444444443
$ws.Range("K125").Value = 665.03506811111106
$ws.Range("L125").Value = 80.288159777777793
$ws.Range("M125").Value = 1407.6282552222222
$ws.Range("N125").Value = 948.91337744444434
$ws.Range("O125").Value = 3991.0678982222216
$ws.Range("P125").Value = 27393.541449666667
$ws.Range("Q125").Value = 2.373425444444444
$ws.Range("R125").Value = 6.5255555555555553E-4

$ws.Range("V125").Value = 2824925103
$ws.Range("W125").Value = 14733975
$ws.Range("X125").Value = 69587931
$ws.Range("Y125").Formula = "=SUM(V125:X125)"
$ws.Range("Z125").Formula = "=(Y125/3307080000)*1000"

$ws.Range("AB125").Value = 2748807187
$ws.Range("AC125").Value = 11645440
$ws.Range("AD125").Value = 69451195
$ws.Range("AE125").Formula = "=SUM(AB125:AD125)"
$ws.Range("AF125").Formula = "=(AE125/3307080000)*1000"

$ws.Range("AH125").Formula = "=Y125-AE125"
$ws.Range("AI125").Formula = "=Z125-AF125"

$ws.Range("AJ125").Value = "???"

# ---------------------------------------------------------------------------
# 5. Final selection state, matching the author's saved view
# ---------------------------------------------------------------------------
$ws.Range("AJ126").Select()
